{"js": "// Rename the \"N 1\" / \"Summary 2\" / \"N 2\" header cells of the summary\n// table to \"Missing 1\" / \"Summary 0\" / \"Missing 0\" respectively.\nconst replacements = [\n  { find: \"N 1\", replace: \"Missing 1\" },\n  { find: \"Summary 2\", replace: \"Summary 0\" },\n  { find: \"N 2\", replace: \"Missing 0\" },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items/text\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the \"N 1\" / \"Summary 2\" / \"N 2\" header cells of the summary\n# table to \"Missing 1\" / \"Summary 0\" / \"Missing 0\" respectively.\n$d = $word.ActiveDocument\n\n$replacements = @{\n    \"N 1\"       = \"Missing 1\"\n    \"Summary 2\" = \"Summary 0\"\n    \"N 2\"       = \"Missing 0\"\n}\n\n$table = $d.Tables.Item(1)\n$headerRow = $table.Rows.Item(1)\n$cellCount = $headerRow.Cells.Count\n\nfor ($i = 1; $i -le $cellCount; $i++) {\n    $cell = $table.Cell(1, $i)\n    # Cell.Range.Text carries a trailing cell-mark (CR + BEL); strip it\n    # before comparing against the plain header labels.\n    $text = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($replacements.ContainsKey($text)) {\n        $cell.Range.Text = $replacements[$text]\n    }\n}\n"}
